$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T050006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2013%3A00%3A06%20GMT&X-Amz-Signature=3aa96b48201e2a565aac3f41d9a703c268f435a10c81a7532b74fc6f6ed17182'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T050006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2013%3A00%3A06%20GMT&X-Amz-Signature=c15f529c02752010a322a3ed58a3b2935c7bc5831bb0f1289609850bc58190fa'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T050006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2013%3A00%3A06%20GMT&X-Amz-Signature=98a6c30aac94f86b5188806d84593e023bf05d43df4c07abe8136c9c7b6f8bb0'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T050006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2013%3A00%3A06%20GMT&X-Amz-Signature=df91a35ebc94579f6770e69e2932829997d1ec8cdf2134ef0a4aac5b2660c98e'
